# 실습과제 : ACME Work Items Hash 데이터 추가
# Adds a "Hash" value for each client row (J2:J21) on the "WI_가공_세부" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WI_가공_세부")

$hashes = @(
    "4882ebf8c03df53c30e512a49d681b132b0a38eb",
    "b2c37847cdd619e4eb57091d1cbd8cac6e9b03a6",
    "eda668adbed17c579e580282b6e5b1a72700acbf",
    "07e7c66179907494b3c4fbcadcc5bfc0c2399550",
    "43b4a48c73a95106f8f5b6eaaeb61b0142b73df2",
    "c641355d0f7fafdd6037fb19303c8363794c6ae3",
    "dafc9516e64e98e9bfd379d9e5b26e77f715cfcb",
    "94ac0eb0dce0fbe7548eafb7a93342ff51470346",
    "7b3377d3426500876741838c26f83b332c54f4c0",
    "f2a6a7a06026ac3d05b6b84c93505b23cb23d303",
    "9b706d312182d6588267c852c9d5ecfc22208d34",
    "88872ad53682de9e00c72c54e72a5ceda130f5d6",
    "30d95d99ad2c100c596bddbac7880773896f82be",
    "1c514c4689f77b487ebcebb3a092bbbff02c7fef",
    "ce2df4044c200f94c8dcdb14e9a1191085fc912a",
    "d408e344ca46072d44f04befc0a9c4efe7dbf256",
    "4b689b8903814a74ccae87c1ae4453ee01560ec1",
    "7cc4773d7fa298cab1f5c637d00bca23e187d67b",
    "7da2444497b885594cd80f066454d5d9e349f78a",
    "96808e6572a5befff0508e654d8b42537b0aaabe"
)

for ($i = 0; $i -lt $hashes.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 10).Value = $hashes[$i]
}
